# Append five new list-paragraph bullet items to the end of the document,
# matching the existing "ListParagraph" / numId=1 bullet list style.

$d = $word.ActiveDocument

$newTexts = @(
    "Updated password confirmation test so it doesn’t redirect to links page if wrong password is given. We need to not redirect the user is wrong password confirmation is given, so adjust controller to display same form again with an error message. ",
    "Installed Sinatra flash into Gemfile. ",
    "Registered Sinatra flash at the top of controller. ",
    "Included flash notice in layout ",
    "Test now passes. "
)

foreach ($text in $newTexts) {
    $last = $d.Paragraphs.Last
    $last.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    $newPara.Range.Text = $text
}
